$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data per commit diff.
# Force text format on numeric-looking Price (D) values so they are stored
# as text strings (matching original inlineStr cell type) rather than numbers.

$ws.Range("D2").Value = "43.845.23"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "2.234.79"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "272.82"
$ws.Range("E5").Value = "  +5.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.69"
$ws.Range("E6").Value = "  +8.99%  "
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.83"
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.74"
$ws.Range("E12").Value = "  +8.38%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "2.573.10"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.89"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "2.236.93"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "43.779.13"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.30"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.07"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.66"
$ws.Range("E24").Value = "  -9.86%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.55"
$ws.Range("E27").Value = "  +5.42%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.78"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.32"
$ws.Range("E29").Value = "  +5.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.07"
$ws.Range("E30").Value = "  -8.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.61"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0903"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0351"
$ws.Range("E37").Value = "  -3.72%  "
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  +17.01%  "
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.38"
$ws.Range("E41").Value = "  -8.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.72"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.47"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0986"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.96"
$ws.Range("E47").Value = "  -5.48%  "
$ws.Range("E48").Value = "  +3.15%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.48"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.426"
$ws.Range("E51").Value = "  -10.88%  "
